$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for columns L (12), M (13), N (14) to width 25
# (ColumnWidth must be set to target-5/6 because Excel stores the OOXML
#  "width" attribute as ColumnWidth + 5/6 for the default Calibri 11 font)
$ws.Range("L1").EntireColumn.ColumnWidth = 24.1666666666667
$ws.Range("M1").EntireColumn.ColumnWidth = 24.1666666666667
$ws.Range("N1").EntireColumn.ColumnWidth = 24.1666666666667

# Force column D to Text format so date-like strings are not auto-converted to date serials
$ws.Range("D2:D46").NumberFormat = "@"

# Row 2 (YT)
$ws.Range("D2").Value = "2026-02-14"
$ws.Range("E2").Value = "2026-02-14 00:18:49"
$ws.Range("F2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YT&dia=2026-02-14T09:00Z"
$ws.Range("G2").Value = "sense dades"
$ws.Range("H2").Value = "sense dades"
$ws.Range("I2").Value = "sense dades"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "sense dades"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "sense dades sense dades"
$ws.Range("N2").Value = "sense dades sense dades"
$ws.Range("O2").Value = "sense dades"

# Row 3 (Z1)
$ws.Range("D3").Value = "2026-02-14"
$ws.Range("E3").Value = "2026-02-14 00:18:51"
$ws.Range("F3").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z1&dia=2026-02-14T09:00Z"
$ws.Range("G3").Value = "sense dades"
$ws.Range("H3").Value = "sense dades"
$ws.Range("I3").Value = "sense dades"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "sense dades"
$ws.Range("L3").Value = "sense dades sense dades"
$ws.Range("M3").Value = "sense dades sense dades"
$ws.Range("N3").Value = "sense dades sense dades"
$ws.Range("O3").Value = "sense dades"

# Row 4 (DN)
$ws.Range("D4").Value = "2026-02-14"
$ws.Range("E4").Value = "2026-02-14 00:18:54"
$ws.Range("F4").Value = "https://www.meteo.cat/observacions/xema/dades?codi=DN&dia=2026-02-14T09:00Z"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "sense dades"
$ws.Range("I4").Value = "sense dades"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "sense dades"
$ws.Range("L4").Value = "sense dades sense dades"
$ws.Range("M4").Value = "sense dades sense dades"
$ws.Range("N4").Value = "sense dades sense dades"
$ws.Range("O4").Value = "sense dades"

# Row 5 (Z6)
$ws.Range("D5").Value = "2026-02-14"
$ws.Range("E5").Value = "2026-02-14 00:18:56"
$ws.Range("F5").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z6&dia=2026-02-14T09:00Z"
$ws.Range("G5").Value = "sense dades"
$ws.Range("H5").Value = "sense dades"
$ws.Range("I5").Value = "sense dades"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "sense dades"
$ws.Range("L5").Value = "sense dades sense dades"
$ws.Range("M5").Value = "sense dades sense dades"
$ws.Range("N5").Value = "sense dades sense dades"
$ws.Range("O5").Value = "sense dades"

# Row 6 (DJ)
$ws.Range("D6").Value = "2026-02-14"
$ws.Range("E6").Value = "2026-02-14 00:18:59"
$ws.Range("F6").Value = "https://www.meteo.cat/observacions/xema/dades?codi=DJ&dia=2026-02-14T09:00Z"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "sense dades"
$ws.Range("I6").Value = "sense dades"
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "sense dades"
$ws.Range("L6").Value = "sense dades sense dades"
$ws.Range("M6").Value = "sense dades sense dades"
$ws.Range("N6").Value = "sense dades sense dades"
$ws.Range("O6").Value = "sense dades"

# Row 7 (X4)
$ws.Range("D7").Value = "2026-02-14"
$ws.Range("E7").Value = "2026-02-14 00:19:01"
$ws.Range("F7").Value = "https://www.meteo.cat/observacions/xema/dades?codi=X4&dia=2026-02-14T09:00Z"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = "sense dades"
$ws.Range("I7").Value = "sense dades"
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = "sense dades"
$ws.Range("L7").Value = "sense dades sense dades"
$ws.Range("M7").Value = "sense dades sense dades"
$ws.Range("N7").Value = "sense dades sense dades"
$ws.Range("O7").Value = "sense dades"

# Row 8 (D5)
$ws.Range("D8").Value = "2026-02-14"
$ws.Range("E8").Value = "2026-02-14 00:19:03"
$ws.Range("F8").Value = "https://www.meteo.cat/observacions/xema/dades?codi=D5&dia=2026-02-14T09:00Z"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "sense dades"
$ws.Range("I8").Value = "sense dades"
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = "sense dades"
$ws.Range("L8").Value = "sense dades sense dades"
$ws.Range("M8").Value = "sense dades sense dades"
$ws.Range("N8").Value = "sense dades sense dades"
$ws.Range("O8").Value = "sense dades"

# Row 9 (YS)
$ws.Range("D9").Value = "2026-02-14"
$ws.Range("E9").Value = "2026-02-14 00:19:06"
$ws.Range("F9").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YS&dia=2026-02-14T09:00Z"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = "sense dades"
$ws.Range("I9").Value = "sense dades"
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = "sense dades"
$ws.Range("L9").Value = "sense dades sense dades"
$ws.Range("M9").Value = "sense dades sense dades"
$ws.Range("N9").Value = "sense dades sense dades"
$ws.Range("O9").Value = "sense dades"

# Row 10 (UN)
$ws.Range("D10").Value = "2026-02-14"
$ws.Range("E10").Value = "2026-02-14 00:19:08"
$ws.Range("F10").Value = "https://www.meteo.cat/observacions/xema/dades?codi=UN&dia=2026-02-14T09:00Z"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "sense dades"
$ws.Range("I10").Value = "sense dades"
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = "sense dades"
$ws.Range("L10").Value = "sense dades sense dades"
$ws.Range("M10").Value = "sense dades sense dades"
$ws.Range("N10").Value = "sense dades sense dades"
$ws.Range("O10").Value = "sense dades"

# Row 11 (MS)
$ws.Range("D11").Value = "2026-02-14"
$ws.Range("E11").Value = "2026-02-14 00:19:11"
$ws.Range("F11").Value = "https://www.meteo.cat/observacions/xema/dades?codi=MS&dia=2026-02-14T09:00Z"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "sense dades"
$ws.Range("I11").Value = "sense dades"
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = "sense dades sense dades"
$ws.Range("N11").Value = "sense dades sense dades"
$ws.Range("O11").Value = "sense dades"

# Row 12 (W1)
$ws.Range("D12").Value = "2026-02-14"
$ws.Range("E12").Value = "2026-02-14 00:19:13"
$ws.Range("F12").Value = "https://www.meteo.cat/observacions/xema/dades?codi=W1&dia=2026-02-14T09:00Z"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = "sense dades"
$ws.Range("I12").Value = "sense dades"
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "sense dades sense dades"
$ws.Range("N12").Value = "sense dades sense dades"
$ws.Range("O12").Value = "sense dades"

# Row 13 (DP)
$ws.Range("D13").Value = "2026-02-14"
$ws.Range("E13").Value = "2026-02-14 00:19:16"
$ws.Range("F13").Value = "https://www.meteo.cat/observacions/xema/dades?codi=DP&dia=2026-02-14T09:00Z"
$ws.Range("G13").Value = "sense dades"
$ws.Range("H13").Value = "sense dades"
$ws.Range("I13").Value = "sense dades"
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = "sense dades"
$ws.Range("L13").Value = "sense dades sense dades"
$ws.Range("M13").Value = "sense dades sense dades"
$ws.Range("N13").Value = "sense dades sense dades"
$ws.Range("O13").Value = "sense dades"

# Row 14 (XL)
$ws.Range("D14").Value = "2026-02-14"
$ws.Range("E14").Value = "2026-02-14 00:19:18"
$ws.Range("F14").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XL&dia=2026-02-14T09:00Z"
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = "sense dades"
$ws.Range("I14").Value = "sense dades"
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = "sense dades"
$ws.Range("L14").Value = "sense dades sense dades"
$ws.Range("M14").Value = "sense dades sense dades"
$ws.Range("N14").Value = "sense dades sense dades"
$ws.Range("O14").Value = "sense dades"

# Row 15 (VZ)
$ws.Range("D15").Value = "2026-02-14"
$ws.Range("E15").Value = "2026-02-14 00:19:21"
$ws.Range("F15").Value = "https://www.meteo.cat/observacions/xema/dades?codi=VZ&dia=2026-02-14T09:00Z"
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = "sense dades"
$ws.Range("I15").Value = "sense dades"
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = "sense dades sense dades"
$ws.Range("N15").Value = "sense dades sense dades"
$ws.Range("O15").Value = "sense dades"

# Row 16 (Z7)
$ws.Range("D16").Value = "2026-02-14"
$ws.Range("E16").Value = "2026-02-14 00:19:23"
$ws.Range("F16").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z7&dia=2026-02-14T09:00Z"
$ws.Range("G16").Value = "sense dades"
$ws.Range("H16").Value = "sense dades"
$ws.Range("I16").Value = "sense dades"
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = "sense dades"
$ws.Range("L16").Value = "sense dades sense dades"
$ws.Range("M16").Value = "sense dades sense dades"
$ws.Range("N16").Value = "sense dades sense dades"
$ws.Range("O16").Value = "sense dades"

# Row 17 (XK)
$ws.Range("D17").Value = "2026-02-14"
$ws.Range("E17").Value = "2026-02-14 00:19:26"
$ws.Range("F17").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XK&dia=2026-02-14T09:00Z"
$ws.Range("G17").Value = "sense dades"
$ws.Range("H17").Value = "sense dades"
$ws.Range("I17").Value = "sense dades"
$ws.Range("J17").Value = ""
$ws.Range("K17").Value = "sense dades"
$ws.Range("L17").Value = "sense dades sense dades"
$ws.Range("M17").Value = "sense dades sense dades"
$ws.Range("N17").Value = "sense dades sense dades"
$ws.Range("O17").Value = "sense dades"

# Row 18 (XJ)
$ws.Range("D18").Value = "2026-02-14"
$ws.Range("E18").Value = "2026-02-14 00:19:28"
$ws.Range("F18").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-14T09:00Z"
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = "sense dades"
$ws.Range("I18").Value = "sense dades"
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = "sense dades"
$ws.Range("L18").Value = "sense dades sense dades"
$ws.Range("M18").Value = "sense dades sense dades"
$ws.Range("N18").Value = "sense dades sense dades"
$ws.Range("O18").Value = "sense dades"

# Row 19 (YU)
$ws.Range("D19").Value = "2026-02-14"
$ws.Range("E19").Value = "2026-02-14 00:19:30"
$ws.Range("F19").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YU&dia=2026-02-14T09:00Z"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = "sense dades"
$ws.Range("I19").Value = "sense dades"
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = "sense dades"
$ws.Range("L19").Value = "sense dades sense dades"
$ws.Range("M19").Value = "sense dades sense dades"
$ws.Range("N19").Value = "sense dades sense dades"
$ws.Range("O19").Value = "sense dades"

# Row 20 (ZE)
$ws.Range("D20").Value = "2026-02-14"
$ws.Range("E20").Value = "2026-02-14 00:19:33"
$ws.Range("F20").Value = "https://www.meteo.cat/observacions/xema/dades?codi=ZE&dia=2026-02-14T09:00Z"
$ws.Range("G20").Value = "sense dades"
$ws.Range("H20").Value = "sense dades"
$ws.Range("I20").Value = "sense dades"
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = "sense dades"
$ws.Range("L20").Value = "sense dades sense dades"
$ws.Range("M20").Value = "sense dades sense dades"
$ws.Range("N20").Value = "sense dades sense dades"
$ws.Range("O20").Value = "sense dades"

# Row 21 (CD)
$ws.Range("D21").Value = "2026-02-14"
$ws.Range("E21").Value = "2026-02-14 00:19:35"
$ws.Range("F21").Value = "https://www.meteo.cat/observacions/xema/dades?codi=CD&dia=2026-02-14T09:00Z"
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = "sense dades"
$ws.Range("I21").Value = "sense dades"
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = "sense dades"
$ws.Range("L21").Value = "sense dades sense dades"
$ws.Range("M21").Value = "sense dades sense dades"
$ws.Range("N21").Value = "sense dades sense dades"
$ws.Range("O21").Value = "sense dades"

# Row 22 (Z2)
$ws.Range("D22").Value = "2026-02-14"
$ws.Range("E22").Value = "2026-02-14 00:19:38"
$ws.Range("F22").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z2&dia=2026-02-14T09:00Z"
$ws.Range("G22").Value = "sense dades"
$ws.Range("H22").Value = "sense dades"
$ws.Range("I22").Value = "sense dades"
$ws.Range("J22").Value = ""
$ws.Range("K22").Value = "sense dades"
$ws.Range("L22").Value = "sense dades sense dades"
$ws.Range("M22").Value = "sense dades sense dades"
$ws.Range("N22").Value = "sense dades sense dades"
$ws.Range("O22").Value = "sense dades"

# Row 23 (Z5)
$ws.Range("D23").Value = "2026-02-14"
$ws.Range("E23").Value = "2026-02-14 00:19:40"
$ws.Range("F23").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z5&dia=2026-02-14T09:00Z"
$ws.Range("G23").Value = "sense dades"
$ws.Range("H23").Value = "sense dades"
$ws.Range("I23").Value = "sense dades"
$ws.Range("J23").Value = ""
$ws.Range("K23").Value = "sense dades"
$ws.Range("L23").Value = "sense dades sense dades"
$ws.Range("M23").Value = "sense dades sense dades"
$ws.Range("N23").Value = "sense dades sense dades"
$ws.Range("O23").Value = "sense dades"

# Row 24 (VK)
$ws.Range("D24").Value = "2026-02-14"
$ws.Range("E24").Value = "2026-02-14 00:19:43"
$ws.Range("F24").Value = "https://www.meteo.cat/observacions/xema/dades?codi=VK&dia=2026-02-14T09:00Z"
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = "sense dades"
$ws.Range("I24").Value = "sense dades"
$ws.Range("J24").Value = ""
$ws.Range("K24").Value = "sense dades"
$ws.Range("L24").Value = "sense dades sense dades"
$ws.Range("M24").Value = "sense dades sense dades"
$ws.Range("N24").Value = "sense dades sense dades"
$ws.Range("O24").Value = "sense dades"

# Row 25 (Z3)
$ws.Range("D25").Value = "2026-02-14"
$ws.Range("E25").Value = "2026-02-14 00:19:45"
$ws.Range("F25").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z3&dia=2026-02-14T09:00Z"
$ws.Range("G25").Value = "sense dades"
$ws.Range("H25").Value = "sense dades"
$ws.Range("I25").Value = "sense dades"
$ws.Range("J25").Value = ""
$ws.Range("K25").Value = "sense dades"
$ws.Range("L25").Value = "sense dades sense dades"
$ws.Range("M25").Value = "sense dades sense dades"
$ws.Range("N25").Value = "sense dades sense dades"
$ws.Range("O25").Value = "sense dades"

# Row 26 (CG)
$ws.Range("D26").Value = "2026-02-14"
$ws.Range("E26").Value = "2026-02-14 00:19:48"
$ws.Range("F26").Value = "https://www.meteo.cat/observacions/xema/dades?codi=CG&dia=2026-02-14T09:00Z"
$ws.Range("G26").Value = "sense dades"
$ws.Range("H26").Value = "sense dades"
$ws.Range("I26").Value = "sense dades"
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = "sense dades"
$ws.Range("L26").Value = "sense dades sense dades"
$ws.Range("M26").Value = "sense dades sense dades"
$ws.Range("N26").Value = "sense dades sense dades"
$ws.Range("O26").Value = "sense dades"

# Row 27 (Z9)
$ws.Range("D27").Value = "2026-02-14"
$ws.Range("E27").Value = "2026-02-14 00:19:50"
$ws.Range("F27").Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z9&dia=2026-02-14T09:00Z"
$ws.Range("G27").Value = "sense dades"
$ws.Range("H27").Value = "sense dades"
$ws.Range("I27").Value = "sense dades"
$ws.Range("J27").Value = ""
$ws.Range("K27").Value = "sense dades"
$ws.Range("L27").Value = "sense dades sense dades"
$ws.Range("M27").Value = "sense dades sense dades"
$ws.Range("N27").Value = "sense dades sense dades"
$ws.Range("O27").Value = "sense dades"

# Row 28 (YB)
$ws.Range("D28").Value = "2026-02-14"
$ws.Range("E28").Value = "2026-02-14 00:19:53"
$ws.Range("F28").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YB&dia=2026-02-14T09:00Z"
$ws.Range("G28").Value = ""
$ws.Range("H28").Value = "sense dades"
$ws.Range("I28").Value = "sense dades"
$ws.Range("J28").Value = ""
$ws.Range("K28").Value = "sense dades"
$ws.Range("L28").Value = "sense dades sense dades"
$ws.Range("M28").Value = "sense dades sense dades"
$ws.Range("N28").Value = "sense dades sense dades"
$ws.Range("O28").Value = "sense dades"

# Row 29 (YP)
$ws.Range("D29").Value = "2026-02-14"
$ws.Range("E29").Value = "2026-02-14 00:19:55"
$ws.Range("F29").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YP&dia=2026-02-14T09:00Z"
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = "sense dades"
$ws.Range("I29").Value = "sense dades"
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = "sense dades"
$ws.Range("L29").Value = "sense dades sense dades"
$ws.Range("M29").Value = "sense dades sense dades"
$ws.Range("N29").Value = "sense dades sense dades"
$ws.Range("O29").Value = "sense dades"

# Row 30 (J5)
$ws.Range("D30").Value = "2026-02-14"
$ws.Range("E30").Value = "2026-02-14 00:19:58"
$ws.Range("F30").Value = "https://www.meteo.cat/observacions/xema/dades?codi=J5&dia=2026-02-14T09:00Z"
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = "sense dades"
$ws.Range("I30").Value = "sense dades"
$ws.Range("J30").Value = ""
$ws.Range("K30").Value = "sense dades"
$ws.Range("L30").Value = "sense dades sense dades"
$ws.Range("M30").Value = "sense dades sense dades"
$ws.Range("N30").Value = "sense dades sense dades"
$ws.Range("O30").Value = "sense dades"

# Row 31 (D6)
$ws.Range("D31").Value = "2026-02-14"
$ws.Range("E31").Value = "2026-02-14 00:20:00"
$ws.Range("F31").Value = "https://www.meteo.cat/observacions/xema/dades?codi=D6&dia=2026-02-14T09:00Z"
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = "sense dades"
$ws.Range("I31").Value = "sense dades"
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = "sense dades"
$ws.Range("L31").Value = "sense dades sense dades"
$ws.Range("M31").Value = "sense dades sense dades"
$ws.Range("N31").Value = "sense dades sense dades"
$ws.Range("O31").Value = "sense dades"

# Row 32 (XR)
$ws.Range("D32").Value = "2026-02-14"
$ws.Range("E32").Value = "2026-02-14 00:20:03"
$ws.Range("F32").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XR&dia=2026-02-14T09:00Z"
$ws.Range("G32").Value = "sense dades"
$ws.Range("H32").Value = "sense dades"
$ws.Range("I32").Value = "sense dades"
$ws.Range("J32").Value = ""
$ws.Range("K32").Value = "sense dades"
$ws.Range("L32").Value = "sense dades sense dades"
$ws.Range("M32").Value = "sense dades sense dades"
$ws.Range("N32").Value = "sense dades sense dades"
$ws.Range("O32").Value = "sense dades"

# Row 33 (YA)
$ws.Range("D33").Value = "2026-02-14"
$ws.Range("E33").Value = "2026-02-14 00:20:05"
$ws.Range("F33").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YA&dia=2026-02-14T09:00Z"
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = "sense dades"
$ws.Range("I33").Value = "sense dades"
$ws.Range("J33").Value = ""
$ws.Range("K33").Value = "sense dades"
$ws.Range("L33").Value = "sense dades sense dades"
$ws.Range("M33").Value = "sense dades sense dades"
$ws.Range("N33").Value = "sense dades sense dades"
$ws.Range("O33").Value = "sense dades"

# Row 34 (DG)
$ws.Range("D34").Value = "2026-02-14"
$ws.Range("E34").Value = "2026-02-14 00:20:08"
$ws.Range("F34").Value = "https://www.meteo.cat/observacions/xema/dades?codi=DG&dia=2026-02-14T09:00Z"
$ws.Range("G34").Value = "sense dades"
$ws.Range("H34").Value = "sense dades"
$ws.Range("I34").Value = "sense dades"
$ws.Range("J34").Value = ""
$ws.Range("K34").Value = "sense dades"
$ws.Range("L34").Value = "sense dades sense dades"
$ws.Range("M34").Value = "sense dades sense dades"
$ws.Range("N34").Value = "sense dades sense dades"
$ws.Range("O34").Value = "sense dades"

# Row 35 (X5)
$ws.Range("D35").Value = "2026-02-14"
$ws.Range("E35").Value = "2026-02-14 00:20:10"
$ws.Range("F35").Value = "https://www.meteo.cat/observacions/xema/dades?codi=X5&dia=2026-02-14T09:00Z"
$ws.Range("G35").Value = "sense dades"
$ws.Range("H35").Value = "sense dades"
$ws.Range("I35").Value = "sense dades"
$ws.Range("J35").Value = ""
$ws.Range("K35").Value = "sense dades"
$ws.Range("L35").Value = "sense dades sense dades"
$ws.Range("M35").Value = "sense dades sense dades"
$ws.Range("N35").Value = "sense dades sense dades"
$ws.Range("O35").Value = "sense dades"

# Row 36 (D4)
$ws.Range("D36").Value = "2026-02-14"
$ws.Range("E36").Value = "2026-02-14 00:20:12"
$ws.Range("F36").Value = "https://www.meteo.cat/observacions/xema/dades?codi=D4&dia=2026-02-14T09:00Z"
$ws.Range("G36").Value = ""
$ws.Range("H36").Value = "sense dades"
$ws.Range("I36").Value = "sense dades"
$ws.Range("J36").Value = ""
$ws.Range("K36").Value = "sense dades"
$ws.Range("L36").Value = "sense dades sense dades"
$ws.Range("M36").Value = "sense dades sense dades"
$ws.Range("N36").Value = "sense dades sense dades"
$ws.Range("O36").Value = "sense dades"

# Row 37 (CI)
$ws.Range("D37").Value = "2026-02-14"
$ws.Range("E37").Value = "2026-02-14 00:20:15"
$ws.Range("F37").Value = "https://www.meteo.cat/observacions/xema/dades?codi=CI&dia=2026-02-14T09:00Z"
$ws.Range("G37").Value = ""
$ws.Range("H37").Value = "sense dades"
$ws.Range("I37").Value = "sense dades"
$ws.Range("J37").Value = ""
$ws.Range("K37").Value = ""
$ws.Range("L37").Value = "sense dades sense dades"
$ws.Range("M37").Value = "sense dades sense dades"
$ws.Range("N37").Value = "sense dades sense dades"
$ws.Range("O37").Value = "sense dades"

# Row 38 (XS)
$ws.Range("D38").Value = "2026-02-14"
$ws.Range("E38").Value = "2026-02-14 00:20:17"
$ws.Range("F38").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XS&dia=2026-02-14T09:00Z"
$ws.Range("G38").Value = ""
$ws.Range("H38").Value = "sense dades"
$ws.Range("I38").Value = "sense dades"
$ws.Range("J38").Value = ""
$ws.Range("K38").Value = "sense dades"
$ws.Range("L38").Value = "sense dades sense dades"
$ws.Range("M38").Value = "sense dades sense dades"
$ws.Range("N38").Value = "sense dades sense dades"
$ws.Range("O38").Value = "sense dades"

# Row 39 (ZC)
$ws.Range("D39").Value = "2026-02-14"
$ws.Range("E39").Value = "2026-02-14 00:20:20"
$ws.Range("F39").Value = "https://www.meteo.cat/observacions/xema/dades?codi=ZC&dia=2026-02-14T09:00Z"
$ws.Range("G39").Value = "sense dades"
$ws.Range("H39").Value = "sense dades"
$ws.Range("I39").Value = "sense dades"
$ws.Range("J39").Value = ""
$ws.Range("K39").Value = "sense dades"
$ws.Range("L39").Value = "sense dades sense dades"
$ws.Range("M39").Value = "sense dades sense dades"
$ws.Range("N39").Value = "sense dades sense dades"
$ws.Range("O39").Value = "sense dades"

# Row 40 (XH)
$ws.Range("D40").Value = "2026-02-14"
$ws.Range("E40").Value = "2026-02-14 00:20:22"
$ws.Range("F40").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XH&dia=2026-02-14T09:00Z"
$ws.Range("G40").Value = "sense dades"
$ws.Range("H40").Value = "sense dades"
$ws.Range("I40").Value = "sense dades"
$ws.Range("J40").Value = ""
$ws.Range("K40").Value = ""
$ws.Range("L40").Value = ""
$ws.Range("M40").Value = "sense dades sense dades"
$ws.Range("N40").Value = "sense dades sense dades"
$ws.Range("O40").Value = "sense dades"

# Row 41 (XE)
$ws.Range("D41").Value = "2026-02-14"
$ws.Range("E41").Value = "2026-02-14 00:20:25"
$ws.Range("F41").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XE&dia=2026-02-14T09:00Z"
$ws.Range("G41").Value = ""
$ws.Range("H41").Value = "sense dades"
$ws.Range("I41").Value = "sense dades"
$ws.Range("J41").Value = ""
$ws.Range("K41").Value = "sense dades"
$ws.Range("L41").Value = "sense dades sense dades"
$ws.Range("M41").Value = "sense dades sense dades"
$ws.Range("N41").Value = "sense dades sense dades"
$ws.Range("O41").Value = "sense dades"

# Row 42 (UE)
$ws.Range("D42").Value = "2026-02-14"
$ws.Range("E42").Value = "2026-02-14 00:20:27"
$ws.Range("F42").Value = "https://www.meteo.cat/observacions/xema/dades?codi=UE&dia=2026-02-14T09:00Z"
$ws.Range("G42").Value = ""
$ws.Range("H42").Value = "sense dades"
$ws.Range("I42").Value = "sense dades"
$ws.Range("J42").Value = ""
$ws.Range("K42").Value = ""
$ws.Range("L42").Value = ""
$ws.Range("M42").Value = "sense dades sense dades"
$ws.Range("N42").Value = "sense dades sense dades"
$ws.Range("O42").Value = "sense dades"

# Row 43 (XO)
$ws.Range("D43").Value = "2026-02-14"
$ws.Range("E43").Value = "2026-02-14 00:20:30"
$ws.Range("F43").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XO&dia=2026-02-14T09:00Z"
$ws.Range("G43").Value = ""
$ws.Range("H43").Value = "sense dades"
$ws.Range("I43").Value = "sense dades"
$ws.Range("J43").Value = ""
$ws.Range("K43").Value = "sense dades"
$ws.Range("L43").Value = "sense dades sense dades"
$ws.Range("M43").Value = "sense dades sense dades"
$ws.Range("N43").Value = "sense dades sense dades"
$ws.Range("O43").Value = "sense dades"

# Row 44 (VS)
$ws.Range("D44").Value = "2026-02-14"
$ws.Range("E44").Value = "2026-02-14 00:20:32"
$ws.Range("F44").Value = "https://www.meteo.cat/observacions/xema/dades?codi=VS&dia=2026-02-14T09:00Z"
$ws.Range("G44").Value = "sense dades"
$ws.Range("H44").Value = "sense dades"
$ws.Range("I44").Value = "sense dades"
$ws.Range("J44").Value = ""
$ws.Range("K44").Value = "sense dades"
$ws.Range("L44").Value = "sense dades sense dades"
$ws.Range("M44").Value = "sense dades sense dades"
$ws.Range("N44").Value = "sense dades sense dades"
$ws.Range("O44").Value = "sense dades"

# Row 45 (YN)
$ws.Range("D45").Value = "2026-02-14"
$ws.Range("E45").Value = "2026-02-14 00:20:34"
$ws.Range("F45").Value = "https://www.meteo.cat/observacions/xema/dades?codi=YN&dia=2026-02-14T09:00Z"
$ws.Range("G45").Value = "sense dades"
$ws.Range("H45").Value = "sense dades"
$ws.Range("I45").Value = "sense dades"
$ws.Range("J45").Value = ""
$ws.Range("K45").Value = "sense dades"
$ws.Range("L45").Value = "sense dades sense dades"
$ws.Range("M45").Value = "sense dades sense dades"
$ws.Range("N45").Value = "sense dades sense dades"
$ws.Range("O45").Value = "sense dades"

# Row 46 (D7)
$ws.Range("D46").Value = "2026-02-14"
$ws.Range("E46").Value = "2026-02-14 00:20:37"
$ws.Range("F46").Value = "https://www.meteo.cat/observacions/xema/dades?codi=D7&dia=2026-02-14T09:00Z"
$ws.Range("G46").Value = ""
$ws.Range("H46").Value = "sense dades"
$ws.Range("I46").Value = "sense dades"
$ws.Range("J46").Value = ""
$ws.Range("K46").Value = "sense dades"
$ws.Range("L46").Value = "sense dades sense dades"
$ws.Range("M46").Value = "sense dades sense dades"
$ws.Range("N46").Value = "sense dades sense dades"
$ws.Range("O46").Value = "sense dades"

